{"js": "// Update the two-digit / one-digit division problems in the worksheet's\n// single table. The table has 20 rows total: every 4th row (0, 4, 8, 12, 16\n// in 0-based indexing) holds the 5 visible \"NN\u00f7N=\" problems for that block,\n// while the 3 rows following each are intentionally blank spacer rows.\n//\n// Some of the new values (e.g. \"16\u00f79=\") coincide with another cell's\n// *original* text elsewhere in the table, so a naive text-wide find/replace\n// could accidentally clobber an already-correct cell. Addressing each cell\n// directly by (row, column) sidesteps that entirely.\n\nconst newValuesByDataRow = [\n  [\"93\u00f73=\", \"65\u00f73=\", \"81\u00f78=\", \"56\u00f74=\", \"12\u00f73=\"],\n  [\"98\u00f79=\", \"10\u00f78=\", \"15\u00f73=\", \"16\u00f79=\", \"75\u00f74=\"],\n  [\"14\u00f72=\", \"43\u00f73=\", \"69\u00f74=\", \"96\u00f75=\", \"28\u00f78=\"],\n  [\"87\u00f72=\", \"89\u00f79=\", \"80\u00f75=\", \"42\u00f79=\", \"85\u00f79=\"],\n  [\"78\u00f72=\", \"56\u00f78=\", \"59\u00f76=\", \"70\u00f73=\", \"87\u00f76=\"],\n];\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let i = 0; i < dataRowIndexes.length; i++) {\n  const rowIndex = dataRowIndexes[i];\n  const newValues = newValuesByDataRow[i];\n  for (let col = 0; col < newValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = newValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit / one-digit division problems in the first table.\n# Each data row of the table (rows 1, 5, 9, 13, 17) holds 5 problems; the\n# other rows are intentionally blank spacer rows. Target each cell directly\n# by (row, column) so the edit is unambiguous even though some of the new\n# values (e.g. \"16\u00f79=\") coincide with other cells' original text.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$t.Cell(1,1).Range.Text = \"93\u00f73=\"\n$t.Cell(1,2).Range.Text = \"65\u00f73=\"\n$t.Cell(1,3).Range.Text = \"81\u00f78=\"\n$t.Cell(1,4).Range.Text = \"56\u00f74=\"\n$t.Cell(1,5).Range.Text = \"12\u00f73=\"\n$t.Cell(5,1).Range.Text = \"98\u00f79=\"\n$t.Cell(5,2).Range.Text = \"10\u00f78=\"\n$t.Cell(5,3).Range.Text = \"15\u00f73=\"\n$t.Cell(5,4).Range.Text = \"16\u00f79=\"\n$t.Cell(5,5).Range.Text = \"75\u00f74=\"\n$t.Cell(9,1).Range.Text = \"14\u00f72=\"\n$t.Cell(9,2).Range.Text = \"43\u00f73=\"\n$t.Cell(9,3).Range.Text = \"69\u00f74=\"\n$t.Cell(9,4).Range.Text = \"96\u00f75=\"\n$t.Cell(9,5).Range.Text = \"28\u00f78=\"\n$t.Cell(13,1).Range.Text = \"87\u00f72=\"\n$t.Cell(13,2).Range.Text = \"89\u00f79=\"\n$t.Cell(13,3).Range.Text = \"80\u00f75=\"\n$t.Cell(13,4).Range.Text = \"42\u00f79=\"\n$t.Cell(13,5).Range.Text = \"85\u00f79=\"\n$t.Cell(17,1).Range.Text = \"78\u00f72=\"\n$t.Cell(17,2).Range.Text = \"56\u00f78=\"\n$t.Cell(17,3).Range.Text = \"59\u00f76=\"\n$t.Cell(17,4).Range.Text = \"70\u00f73=\"\n$t.Cell(17,5).Range.Text = \"87\u00f76=\"\n"}
